$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 829.97015
$ws.Range("J17").Value = 804.4912
$ws.Range("L17").Value = 2413.4736
$ws.Range("N17").Value = -2749.4736
$ws.Range("H32").Value = 2272.3635
$ws.Range("J32").Value = 2749
$ws.Range("L32").Value = 2749
$ws.Range("N32").Value = -3401
$ws.Range("H64").Value = 3736.25
$ws.Range("J64").Value = 3736.25
$ws.Range("L64").Value = 3736.25
$ws.Range("N64").Value = -4232.25
$ws.Range("H67").Value = 3736.25
$ws.Range("J67").Value = 3736.25
$ws.Range("L67").Value = 3736.25
$ws.Range("N67").Value = -5452.25
$ws.Range("H135").Value = 5707.4
$ws.Range("J135").Value = 11818.182
$ws.Range("L135").Value = 106363.638
$ws.Range("N135").Value = -111433.638
$ws.Range("H137").Value = 4791.75
$ws.Range("I137").Value = 2343.4856
$ws.Range("J137").Value = 21929.6
$ws.Range("K137").Value = 7030.4568
$ws.Range("L137").Value = 65788.79999999999
$ws.Range("M137").Value = -4480.4568
$ws.Range("N137").Value = -70888.79999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2688.9092
$ws.Range("I45").Value = 2008.6666
$ws.Range("K45").Value = 2008.6666
$ws.Range("M45").Value = -1631.6666
$ws.Range("H61").Value = 2243.275
$ws.Range("I61").Value = 2085.0264
$ws.Range("K61").Value = 2085.0264
$ws.Range("M61").Value = -1873.0264
$ws.Range("H63").Value = 2770
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2770
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H80").Value = 40100
$ws.Range("J80").Value = 40110
$ws.Range("L80").Value = 40110
$ws.Range("N80").Value = -42106
$ws.Range("H83").Value = 40100
$ws.Range("J83").Value = 40110
$ws.Range("L83").Value = 120330
$ws.Range("N83").Value = -130314
$ws.Range("H136").Value = 2243.275
$ws.Range("I136").Value = 2085.0264
$ws.Range("K136").Value = 6255.0792
$ws.Range("M136").Value = -3705.0792

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 52.5
$ws.Range("I11").Value = 30
$ws.Range("K11").Value = 30
$ws.Range("M11").Value = 110
$ws.Range("H64").Value = 21235.8
$ws.Range("J64").Value = 1916.3334
$ws.Range("L64").Value = 1916.3334
$ws.Range("N64").Value = -2366.3334
$ws.Range("H67").Value = 21235.8
$ws.Range("J67").Value = 1916.3334
$ws.Range("L67").Value = 1916.3334
$ws.Range("N67").Value = -3476.3334
$ws.Range("H135").Value = 69457
$ws.Range("J135").Value = 69457
$ws.Range("L135").Value = 69457
$ws.Range("N135").Value = -79597

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1315.174
$ws.Range("I22").Value = 847.44446
$ws.Range("J22").Value = 2999
$ws.Range("K22").Value = 847.44446
$ws.Range("L22").Value = 2999
$ws.Range("M22").Value = -497.44446
$ws.Range("N22").Value = -3699
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H80").Value = 11500
$ws.Range("J80").Value = 11500
$ws.Range("L80").Value = 11500
$ws.Range("N80").Value = -13746
$ws.Range("H83").Value = 11500
$ws.Range("J83").Value = 11500
$ws.Range("L83").Value = 34500
$ws.Range("N83").Value = -45732
$ws.Range("H134").Value = 2573
$ws.Range("I134").Value = 2243.9
$ws.Range("K134").Value = 6731.700000000001
$ws.Range("M134").Value = -4196.700000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 200997
$ws.Range("I18").Value = 333395
$ws.Range("J18").Value = 2400
$ws.Range("K18").Value = 1000185
$ws.Range("L18").Value = 7200
$ws.Range("M18").Value = -1000016
$ws.Range("N18").Value = -7538
$ws.Range("H41").Value = 75
$ws.Range("I41").Value = 75
$ws.Range("K41").Value = 225
$ws.Range("M41").Value = 113
$ws.Range("H44").Value = 4493.7334
$ws.Range("I44").Value = 151.5
$ws.Range("J44").Value = 5161.769
$ws.Range("K44").Value = 454.5
$ws.Range("L44").Value = 15485.307
$ws.Range("M44").Value = -56.5
$ws.Range("N44").Value = -16281.307
$ws.Range("H57").Value = 171665.83
$ws.Range("I57").Value = 999997
$ws.Range("J57").Value = 5999.6
$ws.Range("K57").Value = 2999991
$ws.Range("L57").Value = 17998.8
$ws.Range("M57").Value = -2999432
$ws.Range("N57").Value = -19116.8
$ws.Range("H114").Value = 4276.3228
$ws.Range("I114").Value = 422.14285
$ws.Range("J114").Value = 5400.4585
$ws.Range("K114").Value = 1266.42855
$ws.Range("L114").Value = 16201.3755
$ws.Range("M114").Value = 1987.57145
$ws.Range("N114").Value = -22709.3755

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11639.694
$ws.Range("I132").Value = 15337.2
$ws.Range("K132").Value = 46011.60000000001
$ws.Range("M132").Value = -43481.60000000001
$ws.Range("H139").Value = 9344.25
$ws.Range("J139").Value = 9125.666999999999
$ws.Range("L139").Value = 9125.666999999999
$ws.Range("N139").Value = -19405.667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2628.8147
$ws.Range("J22").Value = 2493.1765
$ws.Range("L22").Value = 2493.1765
$ws.Range("N22").Value = -3083.1765
$ws.Range("H27").Value = 2628.8147
$ws.Range("J27").Value = 2493.1765
$ws.Range("L27").Value = 2493.1765
$ws.Range("N27").Value = -2707.1765
$ws.Range("H68").Value = 7738.3335
$ws.Range("I68").Value = 7161.7
$ws.Range("K68").Value = 7161.7
$ws.Range("M68").Value = -6412.7
$ws.Range("H71").Value = 7738.3335
$ws.Range("I71").Value = 7161.7
$ws.Range("K71").Value = 35808.5
$ws.Range("M71").Value = -32064.5
$ws.Range("H133").Value = 82616.336
$ws.Range("J133").Value = 88776.5
$ws.Range("L133").Value = 88776.5
$ws.Range("N133").Value = -93836.5
$ws.Range("H136").Value = 2387.6191
$ws.Range("I136").Value = 2257
$ws.Range("K136").Value = 6771
$ws.Range("M136").Value = -4221

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3766.6875
$ws.Range("I132").Value = 3998.3333
$ws.Range("J132").Value = 3713.2307
$ws.Range("K132").Value = 11994.9999
$ws.Range("L132").Value = 11139.6921
$ws.Range("M132").Value = -9464.999899999999
$ws.Range("N132").Value = -16199.6921
$ws.Range("H136").Value = 1262.1177
$ws.Range("I136").Value = 1187.6
$ws.Range("K136").Value = 3562.8
$ws.Range("M136").Value = -1012.8
